$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the Insights bullet paragraph that currently reads:
#   "Overall, in all years, the number of orders cumulatively
#    increases as the year progresses. The majority of orders made
#    increase starting to August."
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*The majority of orders made increase*") {
        $target = $cand
        break
    }
}

$p = $target
$pStart = $p.Range.Start

# ------------------------------------------------------------------
# Rewrite "The majority of" -> "Most" (narrow in-place replace).
# ------------------------------------------------------------------
$full = $p.Range.Text
$idx = $full.IndexOf("The majority of")
$len = "The majority of".Length
$r = $d.Range($pStart + $idx, $pStart + $idx + $len)
$r.Text = "Most"

Write-Output $p.Range.Text

# ------------------------------------------------------------------
# Force the final sentence to be split into the same run layout Word
# produces when a user retypes a word in the middle of an existing
# run: a boundary right before/after "Most", another right before/after
# "to", and a last one before " August." -- all runs keep identical
# Times New Roman character formatting, only the run boundaries differ.
# We do this by touching (and restoring) a character property on each
# slice, which stops the writer from silently re-merging equally
# formatted neighbours.
# ------------------------------------------------------------------
$full = $p.Range.Text
$mostStart = $pStart + $full.IndexOf("Most")
$mostEnd = $mostStart + "Most".Length
$toIdx = $pStart + $full.IndexOf(" starting ") + " starting ".Length
$tailStart = $toIdx + "to".Length
$pEnd = $p.Range.End - 1

$slice1 = $d.Range($pStart, $mostStart)
$slice1.Font.Bold = $true
$slice1.Font.Bold = $false

$slice2 = $d.Range($mostStart, $mostEnd)
$slice2.Font.Bold = $true
$slice2.Font.Bold = $false

$slice3 = $d.Range($mostEnd, $toIdx)
$slice3.Font.Bold = $true
$slice3.Font.Bold = $false

$slice4 = $d.Range($toIdx, $tailStart)
$slice4.Font.Bold = $true
$slice4.Font.Bold = $false

$slice5 = $d.Range($tailStart, $pEnd)
$slice5.Font.Bold = $true
$slice5.Font.Bold = $false

Write-Output $p.Range.Text
